# Auto-generated: apply scheduled market-price refresh to the Leve profit tables.
# For each changed row we update the price/profit columns (H-N) to the new
# market snapshot. Where a column's new value is blank (no market data),
# the cell is cleared instead of being set to 0/blank text.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1766.9323
$ws.Range("I17").Value = 1188.8889
$ws.Range("J17").Value = 1870.98
$ws.Range("K17").Value = 3566.6667
$ws.Range("L17").Value = 5612.940000000001
$ws.Range("M17").Value = -3398.6667
$ws.Range("N17").Value = -5948.940000000001
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").ClearContents()
$ws.Range("H39").Value = 1516.25
$ws.Range("I39").Value = 355
$ws.Range("J39").Value = 5000
$ws.Range("K39").Value = 1065
$ws.Range("L39").Value = 15000
$ws.Range("M39").Value = -769
$ws.Range("N39").Value = -15592
$ws.Range("H98").Value = 1703.2122
$ws.Range("I98").Value = 1703.2122
$ws.Range("K98").Value = 1703.2122
$ws.Range("M98").Value = -205.2121999999999
$ws.Range("H116").Value = 34161.184
$ws.Range("I116").Value = 38752.047
$ws.Range("K116").Value = 38752.047
$ws.Range("M116").Value = -35310.047
$ws.Range("H122").Value = 1703.2122
$ws.Range("I122").Value = 1703.2122
$ws.Range("K122").Value = 5109.6366
$ws.Range("M122").Value = -2659.6366
$ws.Range("H137").Value = 1236823.6
$ws.Range("I137").Value = 1637
$ws.Range("J137").Value = 1821912
$ws.Range("K137").Value = 4911
$ws.Range("L137").Value = 5465736
$ws.Range("M137").Value = -2361
$ws.Range("N137").Value = -5470836
$ws.Range("H138").Value = 2556.6
$ws.Range("I138").Value = 1743.1818
$ws.Range("J138").Value = 3334.652
$ws.Range("K138").Value = 5229.5454
$ws.Range("L138").Value = 10003.956
$ws.Range("M138").Value = -89.54539999999997
$ws.Range("N138").Value = -20283.956

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3550589.8
$ws.Range("I32").Value = 4117414.5
$ws.Range("J32").Value = 18835.54
$ws.Range("K32").Value = 4117414.5
$ws.Range("L32").Value = 18835.54
$ws.Range("M32").Value = -4117127.5
$ws.Range("N32").Value = -19409.54
$ws.Range("H45").Value = 7122.95
$ws.Range("I45").Value = 5398.4287
$ws.Range("K45").Value = 5398.4287
$ws.Range("M45").Value = -5021.4287
$ws.Range("H54").Value = 22499.5
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H101").Value = 87800.5
$ws.Range("J101").Value = 87800.5
$ws.Range("L101").Value = 87800.5
$ws.Range("N101").Value = -94290.5
$ws.Range("H102").Value = 3932.2666
$ws.Range("I102").Value = 3153.3845
$ws.Range("K102").Value = 3153.3845
$ws.Range("M102").Value = -1531.3845
$ws.Range("H112").Value = 37693
$ws.Range("J112").Value = 37693
$ws.Range("L112").Value = 37693
$ws.Range("N112").Value = -40647
$ws.Range("H132").Value = 3540.5625
$ws.Range("I132").Value = 3443.2666
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 10329.7998
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -7799.799800000001
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2096.3333
$ws.Range("J20").Value = 1676.4
$ws.Range("L20").Value = 1676.4
$ws.Range("N20").Value = -2170.4
$ws.Range("H134").Value = 18176.615
$ws.Range("I134").Value = 19344.535
$ws.Range("K134").Value = 58033.605
$ws.Range("M134").Value = -55498.605

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 5093.625
$ws.Range("I12").Value = 2249.8572
$ws.Range("J12").Value = 25000
$ws.Range("K12").Value = 2249.8572
$ws.Range("L12").Value = 25000
$ws.Range("M12").Value = -2079.8572
$ws.Range("N12").Value = -25340
$ws.Range("H58").Value = 2450.7
$ws.Range("I58").Value = 1500.875
$ws.Range("K58").Value = 1500.875
$ws.Range("M58").Value = -1297.875
$ws.Range("H136").Value = 2450.7
$ws.Range("I136").Value = 1500.875
$ws.Range("K136").Value = 4502.625
$ws.Range("M136").Value = -1952.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 550.3333
$ws.Range("I23").Value = 34.2
$ws.Range("J23").Value = 808.4
$ws.Range("K23").Value = 102.6
$ws.Range("L23").Value = 2425.2
$ws.Range("M23").Value = 132.4
$ws.Range("N23").Value = -2895.2
$ws.Range("H26").Value = 175.16667
$ws.Range("H68").Value = 8393
$ws.Range("I68").Value = 1496
$ws.Range("J68").Value = 10274
$ws.Range("K68").Value = 4488
$ws.Range("L68").Value = 30822
$ws.Range("M68").Value = -3677
$ws.Range("N68").Value = -32444
$ws.Range("H71").Value = 8393
$ws.Range("I71").Value = 1496
$ws.Range("J71").Value = 10274
$ws.Range("K71").Value = 13464
$ws.Range("L71").Value = 92466
$ws.Range("M71").Value = -9408
$ws.Range("N71").Value = -100578
$ws.Range("H75").Value = 4224.9375
$ws.Range("J75").Value = 5415.4546
$ws.Range("L75").Value = 16246.3638
$ws.Range("N75").Value = -18242.3638
$ws.Range("H78").Value = 4224.9375
$ws.Range("J78").Value = 5415.4546
$ws.Range("L78").Value = 48739.0914
$ws.Range("N78").Value = -58723.0914
$ws.Range("H92").Value = 1779
$ws.Range("I92").Value = 2123.75
$ws.Range("K92").Value = 6371.25
$ws.Range("M92").Value = -5123.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2203.2974
$ws.Range("I102").Value = 1141.5385
$ws.Range("K102").Value = 1141.5385
$ws.Range("M102").Value = 480.4614999999999
$ws.Range("H122").Value = 4301.3213
$ws.Range("I122").Value = 4714.7393
$ws.Range("K122").Value = 14144.2179
$ws.Range("M122").Value = -11694.2179
$ws.Range("H132").Value = 756709.5600000001
$ws.Range("J132").Value = 4937
$ws.Range("L132").Value = 14811
$ws.Range("N132").Value = -19871
$ws.Range("H134").Value = 45301.273
$ws.Range("J134").Value = 45301.273
$ws.Range("L134").Value = 135903.819
$ws.Range("N134").Value = -140973.819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 20000
$ws.Range("J98").Value = 20000
$ws.Range("L98").Value = 20000
$ws.Range("N98").Value = -25990
$ws.Range("H100").Value = 7059.227
$ws.Range("I100").Value = 2546.8235
$ws.Range("J100").Value = 22401.4
$ws.Range("K100").Value = 2546.8235
$ws.Range("L100").Value = 22401.4
$ws.Range("M100").Value = -2005.8235
$ws.Range("N100").Value = -23483.4
$ws.Range("H110").Value = 70459.86
$ws.Range("J110").Value = 70459.86
$ws.Range("L110").Value = 70459.86
$ws.Range("N110").Value = -78639.86

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2597.8333
$ws.Range("I96").Value = 1196.6666
$ws.Range("K96").Value = 1196.6666
$ws.Range("M96").Value = 176.3334
$ws.Range("H113").Value = 1233.2106
$ws.Range("I113").Value = 642.13336
$ws.Range("J113").Value = 3449.75
$ws.Range("K113").Value = 1926.40008
$ws.Range("L113").Value = 10349.25
$ws.Range("M113").Value = 243.5999199999999
$ws.Range("N113").Value = -14689.25

